# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
# Drop the bold/border header style entirely (cells go back to the default
# "Normal" style) and clear out the stray "Unnamed: 0" label in A1.
$ws.Range("A1").ClearContents()
$ws.Range("A1:S1").Style = "Normal"

# --- Corrected numeric values (pre/post/total fixation metrics) -----------
# Row 3 - Revisit count
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 8

# Row 4 - Fixation count
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 18
$ws.Range("D4").Value = 4
$ws.Range("F4").Value = 2
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 266

# Row 5 - Dwell time (ms)
$ws.Range("B5").Value = 734.04
$ws.Range("C5").Value = 3303.51
$ws.Range("D5").Value = 1017.74
$ws.Range("F5").Value = 667.4
$ws.Range("L5").Value = 2319.05
$ws.Range("M5").Value = 46551.96

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 0.42
$ws.Range("C6").Value = 1.9
$ws.Range("D6").Value = 0.58
$ws.Range("F6").Value = 0.38
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 26.72

# Row 7 - Fixation duration (ms)
$ws.Range("B7").Value = 367.02
$ws.Range("C7").Value = 183.53
$ws.Range("D7").Value = 254.44
$ws.Range("F7").Value = 333.7
$ws.Range("L7").Value = 231.9
$ws.Range("M7").Value = 175.01

# --- "param" (K) and "variabl1" (O) columns were mis-populated; blank them
#     out for rows 3-8 (row 8 itself only needed this cleanup). ------------
$ws.Range("K3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("O6").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("O7").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("O8").ClearContents()

# --- Drop the trailing fully-blank row 11 ----------------------------------
$ws.Rows(11).Delete()
